$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain literal-text numbers in this sheet
# (e.g. "65.580.24", "0.0000278") with no special cell format. Several of
# the refreshed values parse as ordinary decimals, so Excel would silently
# convert them to real numbers on entry. Flip those cells to Text first,
# write the string, then drop the style back to Normal so the saved file
# keeps the original "no explicit format" cell shape.
$textCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D43", "D44", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.864.29"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "3.411.75"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "563.73"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "177.14"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "3.399.01"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").Value = "0.644"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "53.89"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "0.0000280"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "9.27"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "3.937.74"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "18.37"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "3.384.31"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "65.579.73"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "11.90"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "494.10"
$ws.Range("E22").Value = "  +5.20%  "
$ws.Range("D23").Value = "4.96"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "4.14"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "89.24"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").Value = "14.17"
$ws.Range("E26").Value = "  +5.08%  "
$ws.Range("D27").Value = "2.93"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "8.79"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "31.61"
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").Value = "6.59"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "11.55"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "62.81"
$ws.Range("E33").Value = "  +4.90%  "
$ws.Range("D34").Value = "577.63"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "0.109"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "3.62"
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "36.13"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").Value = "0.375"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "3.125.91"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "2.81"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "0.0419"
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "2.45"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "140.67"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "8.49"
$ws.Range("E51").Value = "  +1.18%  "

# Restore default styling on the text-coerced cells (clears the temporary
# Text number format applied above).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
